$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated symbol list (Price + Volume(1h)) per the "Updated symbol list" GitHub
# Actions commit (Mon Jan 30 09:32:51 UTC 2023). Source data refreshes the Price
# (column D) and Volume(1h) percentage (column E) for each coin row.
#
# The sheet stores these as literal text (e.g. "310.49", "-0.41%"), not numbers,
# so every value below is written with a leading apostrophe, which Excel
# interprets as "store verbatim as text" and strips from the stored value.
$updates = [ordered]@{
    'D2' = '310.49'
    'E2' = '-0.41%'
    'D3' = '37.93'
    'E3' = '-3.59%'
    'D4' = '5.081'
    'E4' = '-1.01%'
    'D5' = '0.07941'
    'E5' = '-2.31%'
    'D6' = '1.987'
    'E6' = '0.30%'
    'D7' = '4.447'
    'E7' = '4.97%'
    'D8' = '8.290'
    'E8' = '1.88%'
    'D9' = '3.090'
    'E9' = '-8.04%'
    'D10' = '0.9338'
    'E10' = '0.64%'
    'D11' = '0.1280'
    'E11' = '-8.76%'
    'D12' = '0.1901'
    'E12' = '-1.48%'
    'D13' = '0.08829'
    'E13' = '-2.20%'
    'D14' = '0.03464'
    'E14' = '-1.48%'
    'D15' = '0.09717'
    'E15' = '-0.99%'
    'D16' = '0.001407'
    'E16' = '0.37%'
    'D17' = '0.006419'
    'E17' = '6.77%'
    'D18' = '3.588'
    'E18' = '-2.61%'
    'D19' = '0.3404'
    'D20' = '0.1295'
    'E20' = '-1.26%'
    'E21' = '8.46%'
    'E22' = '4.58%'
    'D23' = '0.04367'
    'E23' = '0.08%'
    'D24' = '0.001241'
    'E24' = '0.94%'
    'D25' = '0.004667'
    'E25' = '-2.82%'
    'E26' = '176.59%'
    'D39' = '0.02186'
    'E39' = '2.87%'
    'D40' = '0.05086'
    'E40' = '-1.94%'
    'D41' = '0.007591'
    'E41' = '2.04%'
    'D42' = '0.009813'
    'E42' = '-0.01%'
    'D43' = '0.1380'
    'E43' = '0.81%'
    'D44' = '0.002043'
    'E44' = '-3.57%'
    'D45' = '0.008864'
    'E45' = '-1.44%'
    'D46' = '0.00006675'
    'E46' = '4.27%'
    'D47' = '0.00000000757'
    'E47' = '1.06%'
    'D48' = '0.003017'
    'E48' = '18.29%'
    'D49' = '0.001212'
    'E49' = '21.18%'
    'D50' = '0.00002121'
    'E50' = '1.06%'
    'D51' = '0.0002020'
    'E51' = '1.06%'
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = "'" + $updates[$cellRef]
}

# Re-normalize style on the touched cells: writing a leading apostrophe makes
# Excel tag the cell with a "quote prefix" style so it keeps displaying/storing
# the value as text; resetting the style back to Normal removes that bookkeeping
# flag while leaving the text value and its Text cell-type untouched.
$styleRanges = @(
    "D2:E2", "D3:E3", "D4:E4", "D5:E5", "D6:E6", "D7:E7", "D8:E8", "D9:E9", "D10:E10", "D11:E11",
    "D12:E12", "D13:E13", "D14:E14", "D15:E15", "D16:E16", "D17:E17", "D18:E18", "D19", "D20:E20", "E21",
    "E22", "D23:E23", "D24:E24", "D25:E25", "E26", "D39:E39", "D40:E40", "D41:E41", "D42:E42", "D43:E43",
    "D44:E44", "D45:E45", "D46:E46", "D47:E47", "D48:E48", "D49:E49", "D50:E50", "D51:E51"
)
foreach ($rangeRef in $styleRanges) {
    $ws.Range($rangeRef).Style = "Normal"
}
